$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("H3").Value = 2.7
$ws.Range("I3").Value = 2.65

# Row 6 updates
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3.25
$ws.Range("T6").Value = 5.7
$ws.Range("U6").Value = 9.75
$ws.Range("W6").Value = 24
$ws.Range("X6").Value = 24
$ws.Range("AB6").Value = 18.5
$ws.Range("AI6").Value = 37
$ws.Range("AJ6").Value = 55
